# Fruta / hortaliza, semanal
# The data rows (2-43) got reshuffled/reordered (e.g. re-sorted on source
# extraction) while every row's payload (Fecha, Calidad, Volumen, Precio
# minimo/maximo/promedio, Origen, Precio $/Kg) travels together as a unit.
# We snapshot the current values first, then write them back out in the
# new row order so we never read an already-overwritten cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 43

# Columns that move together as a block for each record.
$cols = @(4, 9, 10, 11, 12, 13, 15, 16)   # D, I, J, K, L, M, O, P

# Maps each destination row (2..43) to the row it should pull its
# D/I/J/K/L/M/O/P values from, evaluated against the ORIGINAL layout.
$sourceRow = @{
    2 = 11;  3 = 22;  4 = 21;  5 = 38;  6 = 28;  7 = 9;   8 = 18;  9 = 42;
    10 = 30; 11 = 26; 12 = 24; 13 = 23; 14 = 29; 15 = 31; 16 = 41; 17 = 6;
    18 = 25; 19 = 10; 20 = 15; 21 = 16; 22 = 43; 23 = 27; 24 = 8;  25 = 36;
    26 = 5;  27 = 34; 28 = 37; 29 = 12; 30 = 35; 31 = 33; 32 = 40; 33 = 17;
    34 = 7;  35 = 2;  36 = 3;  37 = 32; 38 = 19; 39 = 20; 40 = 13; 41 = 4;
    42 = 39; 43 = 14
}

# Snapshot the original values for every relevant cell before writing
# anything back, since several source rows are read more than once /
# rows are both a source and a destination.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $sourceRow[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $snapshot["$src-$c"]
    }
}
